# Refresh the cryptocurrency market snapshot (prices + 1h volume deltas) and
# fix the TRON / WrappedEther row ordering that had swapped since the last run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "63.932.68"; ForceText = $false },
    @{ Cell = "E2"; Value = "  -3.84%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "3.351.18"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -4.10%  "; ForceText = $false },
    @{ Cell = "D4"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E4"; Value = "  +0.33%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "556.65"; ForceText = $true },
    @{ Cell = "E5"; Value = "  -0.22%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "171.92"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -7.00%  "; ForceText = $false },
    @{ Cell = "D7"; Value = "0.611"; ForceText = $true },
    @{ Cell = "E7"; Value = "  -3.96%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "3.339.11"; ForceText = $false },
    @{ Cell = "E8"; Value = "  -4.21%  "; ForceText = $false },
    @{ Cell = "E9"; Value = "  +0.05%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.618"; ForceText = $true },
    @{ Cell = "E10"; Value = "  -2.44%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "0.150"; ForceText = $true },
    @{ Cell = "E11"; Value = "  -2.85%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "53.68"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -1.39%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "0.0000264"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -2.74%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "8.89"; ForceText = $true },
    @{ Cell = "E14"; Value = "  -4.32%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "3.886.75"; ForceText = $false },
    @{ Cell = "E15"; Value = "  -3.85%  "; ForceText = $false },
    @{ Cell = "B16"; Value = "WrappedEther"; ForceText = $false },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; ForceText = $false },
    @{ Cell = "D16"; Value = "3.361.00"; ForceText = $false },
    @{ Cell = "E16"; Value = "  -3.65%  "; ForceText = $false },
    @{ Cell = "B17"; Value = "TRON"; ForceText = $false },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; ForceText = $false },
    @{ Cell = "D17"; Value = "0.118"; ForceText = $true },
    @{ Cell = "E17"; Value = "  -2.82%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "17.78"; ForceText = $true },
    @{ Cell = "E18"; Value = "  -3.67%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "11.70"; ForceText = $true },
    @{ Cell = "E19"; Value = "  -2.69%  "; ForceText = $false },
    @{ Cell = "D20"; Value = "63.994.74"; ForceText = $false },
    @{ Cell = "E20"; Value = "  -3.66%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "0.974"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -1.97%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "403.36"; ForceText = $true },
    @{ Cell = "E22"; Value = "  -4.43%  "; ForceText = $false },
    @{ Cell = "E23"; Value = "  +0.65%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "4.26"; ForceText = $true },
    @{ Cell = "E24"; Value = "  +2.49%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "13.35"; ForceText = $true },
    @{ Cell = "E25"; Value = "  +9.13%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "82.81"; ForceText = $true },
    @{ Cell = "E26"; Value = "  -3.68%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "10.65"; ForceText = $true },
    @{ Cell = "E27"; Value = "  -2.32%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "2.74"; ForceText = $true },
    @{ Cell = "E28"; Value = "  -5.24%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "8.69"; ForceText = $true },
    @{ Cell = "E29"; Value = "  -4.55%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "29.18"; ForceText = $true },
    @{ Cell = "E30"; Value = "  -3.27%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "6.40"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -2.88%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "581.04"; ForceText = $true },
    @{ Cell = "E32"; Value = "  -7.34%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "11.34"; ForceText = $true },
    @{ Cell = "E33"; Value = "  -3.35%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "0.106"; ForceText = $true },
    @{ Cell = "E34"; Value = "  -3.98%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "57.95"; ForceText = $true },
    @{ Cell = "E35"; Value = "  -3.27%  "; ForceText = $false },
    @{ Cell = "E36"; Value = "  -1.07%  "; ForceText = $false },
    @{ Cell = "E37"; Value = "  -0.26%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "35.64"; ForceText = $true },
    @{ Cell = "E38"; Value = "  -5.48%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "3.44"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -3.34%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "0.0₃0739"; ForceText = $false },
    @{ Cell = "E40"; Value = "  -9.04%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "0.368"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -4.50%  "; ForceText = $false },
    @{ Cell = "D42"; Value = "3.147.20"; ForceText = $false },
    @{ Cell = "E42"; Value = "  +0.68%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E43"; Value = "  +0.33%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "2.84"; ForceText = $true },
    @{ Cell = "E44"; Value = "  -0.11%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "3.23"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -2.64%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "2.46"; ForceText = $true },
    @{ Cell = "E46"; Value = "  -5.51%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "0.0404"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -2.56%  "; ForceText = $false },
    @{ Cell = "E48"; Value = "  -4.26%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "0.128"; ForceText = $true },
    @{ Cell = "E49"; Value = "  -4.04%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "132.88"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -4.45%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "8.08"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -4.93%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # These new values happen to look like plain numbers (e.g. "1.00",
        # "556.65"); the source data is text, so force the text number
        # format before assigning and restore the cell's normal style
        # afterwards so no stray formatting is left behind.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
